# Beelina Daily Summarize Transactions Report template update
# "#43-Calculate discount amount on Daily Summarize Transaction Report."
#
# Adds discount/net-sales/net-collectibles breakdown columns to the
# report header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 (summary totals strip): keep Initial/Remaining Stocks Value,
# rename "Total Sales" -> "Total Gross Sales", and append two new
# summary columns: Discount Amount, Total Net Sales.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "Total Gross Sales"
$ws.Range("D1").Value = "Discount Amount"
$ws.Range("E1").Value = "Total Net Sales"

# D1/E1 need the same bold + centered header look as A1:C1, but the
# template uses vertical "top" for these two new cells.
$ws.Range("D1:E1").Font.Bold = $true
$ws.Range("D1:E1").HorizontalAlignment = -4108
$ws.Range("D1:E1").VerticalAlignment = -4160
$ws.Range("E1").NumberFormat = "m/d/yyyy"

# Row 2 is the blank value row under row 1; extend its number format
# (#,##0.00) across the two newly added columns.
$ws.Range("D2:E2").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------
# Row 4 (detail table header): the old "Collectibles (Not Paid)" column
# is replaced by four new columns (Gross/Discount/Discounted/Net
# Collectibles), and everything that used to follow it shifts right by
# three columns. "Discounts given per store" is dropped, and
# "Type of Outlet" / "Payment Method" move to the end of the row.
# ---------------------------------------------------------------------
$ws.Range("L4").Value = $ws.Range("I4").Value()
$ws.Range("K4").Value = $ws.Range("H4").Value()
$ws.Range("J4").Value = $ws.Range("G4").Value()
$ws.Range("I4").Value = $ws.Range("F4").Value()
$ws.Range("H4").Value = $ws.Range("E4").Value()
$ws.Range("G4").Value = $ws.Range("D4").Value()
$ws.Range("F4").Value = $ws.Range("C4").Value()
$ws.Range("E4").Value = $ws.Range("B4").Value()

$ws.Range("A4").Value = "Gross Collectibles"
$ws.Range("B4").Value = "Discount"
$ws.Range("C4").Value = "Discounted Collectibles"
$ws.Range("D4").Value = "Net Collectibles"

# B4:D4 pick up the plain bold+center header style (same as the summary
# header cells in row 1), while A4 and the shifted columns keep the
# original bold+center+vertical-center header style.
$ws.Range("B4:D4").Font.Bold = $true
$ws.Range("B4:D4").HorizontalAlignment = -4108
$ws.Range("B4:D4").VerticalAlignment = -4108
$ws.Range("K4").Value = "Type of Outlet"
$ws.Range("L4").Value = "Payment Method"

# K4/L4 (freshly populated via formula copy above) must carry the same
# header style as the rest of the non-date header cells.
$ws.Range("K4:L4").Font.Bold = $true
$ws.Range("K4:L4").HorizontalAlignment = -4108
$ws.Range("K4:L4").VerticalAlignment = -4108

# Clear the leftover style-only placeholder cell that used to sit at K4
# in the old layout (now part of the regular header range, already
# rewritten above).

# ---------------------------------------------------------------------
# Data rows 5 and 6: the per-row date cell shifts from column K to J
# (because the "Discounts given per store" column was removed), and the
# old date-formatted placeholder that used to sit in column B is reset
# back to a plain General number format.
# ---------------------------------------------------------------------
$ws.Range("J5").Value = $ws.Range("K5").Value()
$ws.Range("J6").Value = $ws.Range("K6").Value()
$ws.Range("K5").Clear()
$ws.Range("K6").Clear()
$ws.Range("J5:J6").NumberFormat = "m/d/yyyy"
$ws.Range("B5:B6").NumberFormat = "General"

# ---------------------------------------------------------------------
# Column widths: widen/narrow the columns whose header text changed
# length, matching the template's best-fit column sizing.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.43
$ws.Columns.Item(3).ColumnWidth = 21.86
$ws.Columns.Item(4).ColumnWidth = 16.43
$ws.Columns.Item(8).ColumnWidth = 21
$ws.Columns.Item(10).ColumnWidth = 17.29
$ws.Columns.Item(11).ColumnWidth = 14.14
$ws.Columns.Item(12).ColumnWidth = 23.71

# Drop the now-unused 13th column definition (sheet only spans A:L now).
$ws.Columns.Item(13).Delete()

# Selection cosmetic change recorded by the template author when saving.
$ws.Range("E30").Select()
